# Jan 2014 updates as requested
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above the "Annualised Performance" row (row 9) for the
# new "Half Year to 31 Dec 2013" performance data line, and clone the
# formatting of the preceding yearly-data row (row 8) onto it so it keeps
# the same borders / number formats as the other data rows.
$ws.Rows.Item(9).Insert()
$ws.Range("A8:E8").Copy()
$ws.Range("A9:E9").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Rows.Item(9).RowHeight = 18

# Updated annualised performance figures (now on row 10 after the insert).
$ws.Range("C10").Value2 = "19.3% pa"
$ws.Range("D10").Value2 = "5.1% pa"
$ws.Range("E10").Value2 = "14.2% pa"

# New half-year data row.
$ws.Range("A9").Value2 = "Half Year to"
$ws.Range("B9").Value2 = 41639
$ws.Range("C9").Value2 = 0.25
$ws.Range("D9").Value2 = 0.14
$ws.Range("E9").Value2 = 0.11

# Updated cumulative performance figures (now on row 12 after the insert).
$ws.Range("C12").Value2 = 2.76
$ws.Range("D12").Value2 = 0.45
$ws.Range("E12").Value2 = 2.31

# Move the active selection the same way the source workbook shows.
$ws.Range("C10").Select()
